$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the "8.71" width formatting from columns B:J out to columns B:N ---
$ws.Columns("K:N").ColumnWidth = 7.83

# --- Add the new "2023" column (K) by copying the format of the current last
#     data column (J) for each row, then widening the right edge of the table
#     with a new border now that K is the rightmost column. ---

# Row 3 (year header)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2023
$ws.Range("K3").Borders.Item(10).LineStyle = 1

# Row 4 (Number of employees)
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 1741
$ws.Range("K4").Borders.Item(10).LineStyle = 1

# Row 5 (Of which: Women)
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 538
$ws.Range("K5").Borders.Item(10).LineStyle = 1

# Row 6 (Men)
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 1203
$ws.Range("K6").Borders.Item(10).LineStyle = 1
